# Fix overtime/pay calculation on the "Employees" sheet.
# Overtime hours recorded on the "Work" sheet were not being folded into
# the employees' totals, and the hourly pay rate used for the payout was
# off. This recomputes the OverTime (C) and Pay (D) columns:
#   - OverTime = previous OverTime + sum of OverTime entries for that
#                EmpID from the "Work" sheet
#   - Pay      = previous OverTime * 30 + (sum of OverTime entries from
#                the "Work" sheet) * 15

$wb = $excel.ActiveWorkbook
$wsEmployees = $wb.Worksheets.Item("Employees")
$wsWork = $wb.Worksheets.Item("Work")

# Figure out how many data rows are on the Work sheet.
$workLastRow = $wsWork.Cells.Item($wsWork.Rows.Count, 1).End(-4162).Row

function Get-OvertimeSum($empId) {
    $sum = 0
    for ($r = 2; $r -le $workLastRow; $r++) {
        $idVal = $wsWork.Cells.Item($r, 1).Value2
        if ($idVal -ne $null -and ("$idVal" -eq "$empId")) {
            $otVal = $wsWork.Cells.Item($r, 2).Value2
            if ($otVal -ne $null) {
                $sum += $otVal
            }
        }
    }
    return $sum
}

$empLastRow = $wsEmployees.Cells.Item($wsEmployees.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $empLastRow; $row++) {
    $empId = $wsEmployees.Cells.Item($row, 1).Value2
    if ($empId -eq $null -or "$empId" -eq "") {
        continue
    }

    $oldOvertime = $wsEmployees.Cells.Item($row, 3).Value2
    $workOvertime = Get-OvertimeSum $empId

    $newOvertime = $oldOvertime + $workOvertime
    $newPay = ($oldOvertime * 30) + ($workOvertime * 15)

    $wsEmployees.Cells.Item($row, 3).Value = $newOvertime
    $wsEmployees.Cells.Item($row, 4).Value = $newPay
}
